# Update odds data for "Poland Ekstraklasa" sheet (league bases update 02-03-2024 08:34)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poland Ekstraklasa")

# Row 198 (match id 196): add result columns FTHG/FTAG/FTR and update odds
$ws.Range("H198").Value = 0
$ws.Range("I198").Value = 1
$ws.Range("J198").Value = "A"
$ws.Range("N198").Value = 1.7
$ws.Range("O198").Value = 3.25
$ws.Range("P198").Value = 5
$ws.Range("R198").Value = 2
$ws.Range("S198").Value = 1.85
$ws.Range("T198").Value = 2
$ws.Range("U198").Value = 1.9
$ws.Range("V198").Value = 1.95
$ws.Range("W198").Value = -1
$ws.Range("X198").Value = -1
$ws.Range("Y198").Value = 4
$ws.Range("Z198").Value = -1
$ws.Range("AA198").Value = 0.8500000000000001
$ws.Range("AB198").Value = -1
$ws.Range("AC198").Value = 0.95

# Row 199 (match id 197): add result columns FTHG/FTAG/FTR and update odds
$ws.Range("H199").Value = 3
$ws.Range("I199").Value = 0
$ws.Range("J199").Value = "H"
$ws.Range("N199").Value = 3.3
$ws.Range("O199").Value = 2.875
$ws.Range("P199").Value = 2.25
$ws.Range("R199").Value = 1.85
$ws.Range("S199").Value = 2
$ws.Range("U199").Value = 2
$ws.Range("V199").Value = 1.85
$ws.Range("W199").Value = 2.3
$ws.Range("X199").Value = -1
$ws.Range("Y199").Value = -1
$ws.Range("Z199").Value = 0.8500000000000001
$ws.Range("AA199").Value = -1
$ws.Range("AB199").Value = 1
$ws.Range("AC199").Value = -1

# Row 200 (match id 198): update odds
$ws.Range("N200").Value = 2.5
$ws.Range("P200").Value = 2.45
$ws.Range("R200").Value = 1.925
$ws.Range("S200").Value = 1.925
$ws.Range("T200").Value = 2.75
$ws.Range("U200").Value = 1.975
$ws.Range("V200").Value = 1.875

# Row 201 (match id 199): update odds
$ws.Range("O201").Value = 3
$ws.Range("P201").Value = 3.6
$ws.Range("R201").Value = 1.825
$ws.Range("S201").Value = 2.025
$ws.Range("T201").Value = 2
$ws.Range("U201").Value = 1.875
$ws.Range("V201").Value = 1.975

# Row 202 (match id 200): update odds
$ws.Range("N202").Value = 1.85
$ws.Range("O202").Value = 3.6
$ws.Range("P202").Value = 3.75
$ws.Range("R202").Value = 1.85
$ws.Range("S202").Value = 2
$ws.Range("U202").Value = 2.025
$ws.Range("V202").Value = 1.825

# Row 203 (match id 201): update odds
$ws.Range("N203").Value = 2.3
$ws.Range("P203").Value = 2.7
$ws.Range("Q203").Value = -0.25
$ws.Range("R203").Value = 2.1
$ws.Range("S203").Value = 1.775

# Row 204 (match id 202): update odds
$ws.Range("N204").Value = 2
$ws.Range("P204").Value = 3.3
$ws.Range("R204").Value = 1.8
$ws.Range("S204").Value = 2.05

# Row 205 (match id 203): update odds
$ws.Range("R205").Value = 2.025
$ws.Range("S205").Value = 1.825
